$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU 16 (row 35) and CU 17 (row 36) tasks are now complete ("Hecho").
$ws.Range("F35").Value = "Hecho"
$ws.Range("G35").Value = 1
$ws.Range("AF35").Value = 1

$ws.Range("F36").Value = "Hecho"
$ws.Range("G36").Value = 1
$ws.Range("AF36").Value = 1
